$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "employees": replace the shared formula (48*4) in B2:B6 with
# the literal value 47, then re-fit the columns and update selection.
# ------------------------------------------------------------------
$wsEmployees = $wb.Worksheets.Item("employees")
$wsEmployees.Activate()

$wsEmployees.Range("B2:B6").Value = 47

$wsEmployees.Columns.Item(1).AutoFit() | Out-Null
$wsEmployees.Columns.Item(2).AutoFit() | Out-Null

$wsEmployees.Range("B2:B6").Select()

# ------------------------------------------------------------------
# Sheet "absences": autofit column A and move the selection to C2.
# ------------------------------------------------------------------
$wsAbsences = $wb.Worksheets.Item("absences")
$wsAbsences.Activate()

$wsAbsences.Columns.Item(1).AutoFit() | Out-Null

$wsAbsences.Range("C2").Select()

# ------------------------------------------------------------------
# Sheet "workplace": move the selection to B2:B3.
# ------------------------------------------------------------------
$wsWorkplace = $wb.Worksheets.Item("workplace")
$wsWorkplace.Activate()

$wsWorkplace.Range("B2:B3").Select()

# ------------------------------------------------------------------
# Sheet "requirements": widen the selection to A2:D5. Keep this sheet
# active/activated last so it remains the visible tab on save.
# ------------------------------------------------------------------
$wsRequirements = $wb.Worksheets.Item("requirements")
$wsRequirements.Activate()

$wsRequirements.Range("A2:D5").Select()
